$wb = $excel.ActiveWorkbook

# --- Update the "Okres" text on the "Opis parametrów" sheet ---
$wsOpis = $wb.Worksheets.Item("Opis parametrów")
$wsOpis.Range("A2").Value = "Okres: 26.01.2026 (pon.) - 01.02.2026 (niedz.)"

# --- Add a new transfer row to the "Oddziały" sheet ---
$wsOdd = $wb.Worksheets.Item("Oddziały")

# Column A and B should share the same width (target: 31.7142857142857 chars,
# i.e. 222px at the engine's default-font pixel grid). The COM ColumnWidth
# setter here round-trips through a coarser 1/6-character save grid, so
# 30.8 is the raw input that lands closest (31.666666666666668) to the
# target after serialization.
$wsOdd.Columns.Item(2).ColumnWidth = 30.8

$wsOdd.Range("A5").Value = "27.01.2026, 13, 18:15-19:00, sala: 34"
$wsOdd.Range("B5").Value = "27.01.2026, 11, 16:35-17:20, sala: 34"
$wsOdd.Range("C5").Value = "Nowak Magdalena"
$wsOdd.Range("D5").Value = "-"
$wsOdd.Range("E5").Value = "4B"
$wsOdd.Range("F5").Value = "Matematyka"
$wsOdd.Range("G5").Value = ""
